$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where "TIME IN" (column C) now also gets copied into "TIME OUT" (column D),
# and the whole row is highlighted with the tan/orange fill (matching the
# pre-existing "s=7" cell style already used elsewhere in the sheet).
$rows = @(7, 8, 9, 10, 14, 15, 16, 17, 18)

foreach ($r in $rows) {
    $timeIn = $ws.Range("C$r").Value()
    $ws.Range("D$r").Value = $timeIn
    $ws.Range("A$r`:P$r").Interior.Color = 6737151
}
